$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Release" column (E), header in row 1 ---
$ws.Range("E1").Value = "Release"
$ws.Range("E1").Font.Bold = $true

# --- New log entry, row 10 (2015-10-08) ---
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.4a"
$ws.Range("C10").Value = 5

# --- Row 6 (2015-10-04): shortened description, version bumps recorded in E/F ---
$ws.Range("D6").Value = "JSON parsing. Dynamic views for Meetings."

# --- Finish the new row 10 entry ---
$ws.Range("D10").Value = "Web view & UI tweaks. Video button for testing."

# --- Row 6 release tags ---
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.1"
$ws.Range("F6").Value = "1.2a"

# --- Row 8 (2015-10-06) release tag ---
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.3a"

# --- Apply the same text format to the rest of the Release column ---
$ws.Range("E2:E5").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E11:E25").NumberFormat = "@"

# --- Column widths: Description grew wider, Release column added ---
$ws.Columns("D").ColumnWidth = 55.666666666666664
$ws.Columns("E").ColumnWidth = 7.5

# --- Selection moved by the editor ---
$ws.Range("C11").Select() | Out-Null

Write-Host "Working hours log updated."
